$d = $word.ActiveDocument

# --- Locate the "FAQs" paragraph / run -------------------------------------
$para = $d.Paragraphs(1).Range
$insertPos = $para.End - 1   # position right after "FAQs", before the paragraph mark

# --- Remove the existing _GoBack bookmark (currently wraps the start of the
#     paragraph); we'll recreate it at the new edit point afterwards, which
#     is how Word naturally tracks "_GoBack" after an edit.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# NOTE: adding a *collapsed* bookmark exactly at the end of a paragraph's
# text (i.e. touching the paragraph mark) resolves to the wrong offset in
# this runtime, so a temporary placeholder character is typed after "FAQs"
# first. That keeps the later bookmark-insertion point safely away from the
# paragraph boundary; the placeholder is stripped out again at the end.
$placeholder = $d.Range($insertPos, $insertPos)
$placeholder.InsertAfter("Z")

# --- Type " & Data Sources" right after "FAQs" (before the placeholder) ----
$newTextStart = $insertPos
$d.Range($newTextStart, $newTextStart).InsertAfter(" & Data Sources")

# New text currently inherits a run boundary from the surrounding text; make
# sure it carries the same explicit formatting as the "FAQs" run so it is
# emitted as its own <w:r> with matching <w:rPr>.
$newTextLen = " & Data Sources".Length
$newRunRange = $d.Range($newTextStart, $newTextStart + $newTextLen)
$newRunRange.Font.Name = "Avenir Light"
$newRunRange.Font.NameBi = "Arial"
$newRunRange.Font.Bold = $true

# --- Recreate the _GoBack bookmark, collapsed, right before the placeholder
$bmPos = $d.Content.End - 2   # one before the placeholder "Z", one before the paragraph mark
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Remove the temporary placeholder character -----------------------------
$placeholderRange = $d.Range($d.Content.End - 2, $d.Content.End - 1)
$placeholderRange.Delete()
